# Refresh the cryptos list (prices / 1h volume %) with the latest scrape.
# Price (column D) and Volume(1h) (column E) are stored as text, not numbers,
# so plain-looking numeric strings are written with a leading quote
# (Range.Formula = "'<text>") to force Excel to keep them as text instead of
# silently parsing them into floating point numbers (which would lose
# trailing zeros / exact formatting). Values that already contain more than
# one '.' (e.g. "42.431.05") or the percent strings (which are padded with
# spaces) are never re-interpreted as numbers, so a plain .Value assignment
# is used for those.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.431.05"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "2.511.99"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Formula = "'307.49"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Formula = "'96.66"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Formula = "'0.587"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Formula = "'36.69"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").Formula = "'0.0814"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Formula = "'7.57"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("E13").Value = "  -4.06%  "
$ws.Range("D14").Value = "2.898.34"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Formula = "'15.73"
$ws.Range("E15").Value = "  +8.75%  "
$ws.Range("D16").Value = "2.502.35"
$ws.Range("E16").Value = "  -3.30%  "
$ws.Range("D17").Formula = "'0.862"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "42.460.88"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Formula = "'12.94"
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Formula = "'6.45"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Formula = "'71.50"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Formula = "'254.02"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").Formula = "'2.93"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("D26").Formula = "'26.95"
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +10.63%  "
$ws.Range("D29").Formula = "'10.16"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").Formula = "'37.49"
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("D31").Formula = "'5.95"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").Formula = "'154.15"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Formula = "'19.25"
$ws.Range("E33").Value = "  +5.84%  "
$ws.Range("D34").Formula = "'3.28"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").Formula = "'0.0788"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  -3.81%  "
$ws.Range("E37").Value = "  -5.36%  "
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").Formula = "'0.120"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Formula = "'24.09"
$ws.Range("E40").Value = "  -11.20%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Formula = "'3.41"
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Formula = "'3.88"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").Formula = "'2.04"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").Formula = "'0.998"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Formula = "'0.0302"
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").Value = "2.033.27"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Formula = "'84.59"
$ws.Range("E47").Value = "  -4.57%  "
$ws.Range("D48").Formula = "'8.98"
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("D49").Value = "2.745.37"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("D50").Formula = "'72.98"
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("E51").Value = "  +0.14%  "
